# Edit Card16 sheet: rename "Correction " header to "Correction" and
# add a new "Serviced by " column (O) to the right of it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card16")

# --- Header row ---
# N1: remove trailing space from "Correction "
$ws.Cells.Item(1, 14).Value = "Correction"

# O1: new header "Serviced by " (note trailing space), same style as N1
$ws.Cells.Item(1, 15).Value = "Serviced by "
$ws.Cells.Item(1, 14).Copy()
$ws.Cells.Item(1, 15).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows 2-12 ---
for ($row = 2; $row -le 12; $row++) {
    # N column: previously blank, now filled with "nan" like the other data columns
    $ws.Cells.Item($row, 14).Value = "nan"

    # O column: new blank cell (exists in the sheet but carries no value),
    # mirrored with a no-op formatting touch so Excel materializes the cell.
    $ws.Cells.Item($row, 15).Font.Bold = $false
}

Write-Output "Card16 updated: N1/O1 headers set, N2:N12 filled with 'nan', O2:O12 created blank."
